# Weekly data refresh: insert a new price record at the top of the data
# block (row 12), pushing the existing rows down by one. This mirrors the
# upstream "Fruta / hortaliza, semanal" update that prepends the latest
# weekly observation to the Frambuesa (raspberry) price series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 12; everything below shifts
# down by one (old row 12 -> 13, ..., old row 117 -> 118).
$ws.Rows(12).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44532
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101004
$ws.Range("J12").Value = "Frambuesa"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("Q12").Value = "`$/bandeja 2 kilos"
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 5000
$ws.Range("T12").Value = 2
